# Populate the "Translations - Common" sheet (tags are keyworded too)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")

$ws.Cells.Item(94, 1).Value = 'cs'
$ws.Cells.Item(94, 2).Value = 'market.aroma.create.title'
$ws.Cells.Item(94, 3).Value = 'Nové aroma'
$ws.Cells.Item(95, 1).Value = 'cs'
$ws.Cells.Item(95, 2).Value = '/api/aroma/create.403.title'
$ws.Cells.Item(95, 3).Value = 'K této funkci nemáte přístup.'
$ws.Cells.Item(96, 1).Value = 'cs'
$ws.Cells.Item(96, 2).Value = '/api/aroma/create.403.subtitle'
$ws.Cells.Item(96, 3).Value = 'Pro vytvoření nového aromatu je třeba oprávnění, které bohužel nemáte.'
$ws.Cells.Item(97, 1).Value = 'cs'
$ws.Cells.Item(97, 2).Value = 'shared.certificate.link.button'
$ws.Cells.Item(97, 3).Value = 'Certifikáty'
$ws.Cells.Item(98, 1).Value = 'cs'
$ws.Cells.Item(98, 2).Value = 'shared.license.link.button'
$ws.Cells.Item(98, 3).Value = 'Licence'
$ws.Cells.Item(99, 1).Value = 'cs'
$ws.Cells.Item(99, 2).Value = '/api/aroma/create.name.label'
$ws.Cells.Item(99, 3).Value = 'Název'
$ws.Cells.Item(100, 1).Value = 'cs'
$ws.Cells.Item(100, 2).Value = '/api/aroma/create.name.label.tooltip'
$ws.Cells.Item(100, 3).Value = 'Název aromatu by měl co nejpřesněji odpovídat názvu na lahvičce. Běžně se může stát, že se zde bude také nacházet název edice, např. Moments od Dinner Lady.'
$ws.Cells.Item(101, 1).Value = 'cs'
$ws.Cells.Item(101, 2).Value = '/api/aroma/create.code.label'
$ws.Cells.Item(101, 3).Value = 'Kód'
$ws.Cells.Item(102, 1).Value = 'cs'
$ws.Cells.Item(102, 2).Value = '/api/aroma/create.code.label.tooltip'
$ws.Cells.Item(102, 3).Value = 'Kód aromatu není třeba vyplňovat, dokud si nepřejete uvést nějaký vlastní, konkrétní. Systém jej vygeneruje na pozadí sám.'
$ws.Cells.Item(103, 1).Value = 'cs'
$ws.Cells.Item(103, 2).Value = '/api/aroma/create.create'
$ws.Cells.Item(103, 3).Value = 'Uložit aroma'
$ws.Cells.Item(104, 1).Value = 'cs'
$ws.Cells.Item(104, 2).Value = '/api/aroma/create.vendorId.label'
$ws.Cells.Item(104, 3).Value = 'Výrobce'
$ws.Cells.Item(105, 1).Value = 'cs'
$ws.Cells.Item(105, 2).Value = 'shared.vendor.create.button'
$ws.Cells.Item(105, 3).Value = 'Přidat výrobce'
$ws.Cells.Item(106, 1).Value = 'cs'
$ws.Cells.Item(106, 2).Value = 'shared.vendor.create.title'
$ws.Cells.Item(106, 3).Value = 'Nový výrobce'
$ws.Cells.Item(107, 1).Value = 'cs'
$ws.Cells.Item(107, 2).Value = '/api/vendor/create.name.label'
$ws.Cells.Item(107, 3).Value = 'Jméno'
$ws.Cells.Item(108, 1).Value = 'cs'
$ws.Cells.Item(108, 2).Value = '/api/vendor/create.name.label.tooltip'
$ws.Cells.Item(108, 3).Value = 'Uveďte prosím co nejpřesněji jméno výrobce, např. AEON nebo Dinner Lady.'
$ws.Cells.Item(109, 1).Value = 'cs'
$ws.Cells.Item(109, 2).Value = '/api/vendor/create.create'
$ws.Cells.Item(109, 3).Value = 'Uložit výrobce'
$ws.Cells.Item(110, 1).Value = 'cs'
$ws.Cells.Item(110, 2).Value = '/api/vendor/create.success'
$ws.Cells.Item(110, 3).Value = 'Výrobce [{{name}}] byl úspěšně vytvořen.'
$ws.Cells.Item(111, 1).Value = 'cs'
$ws.Cells.Item(111, 2).Value = '/api/aroma/create.name.label.required'
$ws.Cells.Item(111, 3).Value = 'Název aromatu je povinná položka, bez toho to nejde.'
$ws.Cells.Item(112, 1).Value = 'cs'
$ws.Cells.Item(112, 2).Value = '/api/aroma/create.vendorId.label.required'
$ws.Cells.Item(112, 3).Value = 'Vyberte nebo vytvořte prosím výrobce.'
$ws.Cells.Item(113, 1).Value = 'cs'
$ws.Cells.Item(113, 2).Value = '/api/vendor/create.name.label.required'
$ws.Cells.Item(113, 3).Value = 'Vyplňte prosím jméno výrobce.'
$ws.Cells.Item(114, 1).Value = 'cs'
$ws.Cells.Item(114, 2).Value = '/api/aroma/create.tasteIds.label'
$ws.Cells.Item(114, 3).Value = 'Příchutě'
$ws.Cells.Item(115, 1).Value = 'cs'
$ws.Cells.Item(115, 2).Value = '/api/aroma/create.tasteIds.label.tooltip'
$ws.Cells.Item(115, 3).Value = 'Vyberte hlavní složky příchutě aromatu. Pro pokročilejší hodnocení je možné přidat i položku throat hit pro možnost hodnocení nikotinového kopance.'
$ws.Cells.Item(116, 1).Value = 'cs'
$ws.Cells.Item(116, 2).Value = 'shared.tag.taste.create.button'
$ws.Cells.Item(116, 3).Value = 'Přidat příchuť'
$ws.Cells.Item(117, 1).Value = 'cs'
$ws.Cells.Item(117, 2).Value = 'shared.tag.taste.create.title'
$ws.Cells.Item(118, 1).Value = 'cs'
$ws.Cells.Item(118, 2).Value = '/api/tag/create.403.subtitle'
$ws.Cells.Item(118, 3).Value = 'Pro přidání příchutě bohužel nemáte oprávnění.'
$ws.Cells.Item(117, 3).Value = 'Nová příchuť'
$ws.Cells.Item(119, 1).Value = 'cs'
$ws.Cells.Item(119, 2).Value = '/api/tag/create.403.title'
$ws.Cells.Item(119, 3).Value = 'K této funkci nemáte přístup.'
$ws.Cells.Item(120, 1).Value = 'cs'
$ws.Cells.Item(120, 2).Value = '/api/tag/create.tag.label'
$ws.Cells.Item(120, 3).Value = 'Název'
$ws.Cells.Item(121, 1).Value = 'cs'
$ws.Cells.Item(121, 2).Value = '/api/tag/create.tag.label.tooltip'
$ws.Cells.Item(121, 3).Value = 'Název by měl být malými písmeny, ideálně jednoslovný, případně oddělen pomlčkami (např. foo-bar); používejte prosím anglické výrazy, štítky je možné později přeložit.'
$ws.Cells.Item(122, 1).Value = 'cs'
$ws.Cells.Item(122, 2).Value = '/api/tag/create.tag.label.required'
$ws.Cells.Item(122, 3).Value = 'Název je povinná položka.'
$ws.Cells.Item(123, 1).Value = 'cs'
$ws.Cells.Item(123, 2).Value = '/api/tag/create.sort.label'
$ws.Cells.Item(123, 3).Value = 'Řazení'
$ws.Cells.Item(124, 1).Value = 'cs'
$ws.Cells.Item(124, 2).Value = '/api/tag/create.sort.label.tooltip'
$ws.Cells.Item(124, 3).Value = 'Vyplněním této položky můžete určit pořadí/řazení vybraného štítku mezi ostatními.'
$ws.Cells.Item(125, 1).Value = 'cs'
$ws.Cells.Item(125, 2).Value = '/api/tag/create.group.label'
$ws.Cells.Item(125, 3).Value = 'Skupina'
$ws.Cells.Item(126, 1).Value = 'cs'
$ws.Cells.Item(126, 2).Value = '/api/tag/create.group.label.tooltip'
$ws.Cells.Item(126, 3).Value = 'Vyplňte prosím skupinu, do které štítek patří; aplikace obecně skupiny používá v různých situacích, např. příchutě, typy potahů apod.'
$ws.Cells.Item(127, 1).Value = 'cs'
$ws.Cells.Item(127, 2).Value = '/api/tag/create.group.label.required'
$ws.Cells.Item(127, 3).Value = 'Skupina je povinná položka.'
$ws.Cells.Item(128, 1).Value = 'cs'
$ws.Cells.Item(128, 2).Value = '/api/tag/create.create'
$ws.Cells.Item(128, 3).Value = 'Uložit štítek'
$ws.Cells.Item(129, 1).Value = 'cs'
$ws.Cells.Item(129, 2).Value = '/api/tag/create.shared.tag.create.success'
$ws.Cells.Item(129, 3).Value = 'Štítek [{{tag}}] byl úspěšně vytvořen.'

# Update the active selection / active sheet view to match the final state
$ws.Activate()
$ws.Range("B119").Select() | Out-Null
